# Konnect Bill Payment Verification Checks added
# Adds 6 new "query" columns (W:AB) with header/value pairs to every row
# of the SendMoney.xlsx automation-data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new shared-string values -------------------------------------------------
# (kept in the exact order they first appear so the generated sharedStrings
#  table lines up with the authored workbook: W1, W2, X1, X2, Y2, Z1, Y1,
#  AA1, AA2, AB1, AB2, Z2)
$branchCodeHeader   = "branch_code_query"
$branchCodeSql      = "SELECT DT.BRANCH_CODE FROM DC_TRANSACTION DT WHERE DT.TRANSACTION_ID = '{TRANSACTION_ID}'"
$toBranchHeader     = "to_branch_query"
$toBranchSql        = "SELECT DT.FT_TO_ACCOUNT_BRANCH_CODE FROM DC_TRANSACTION DT WHERE DT.TRANSACTION_ID = '{TRANSACTION_ID}'"
$beneIdSql          = "SELECT BENEFICIARY_ID FROM DC_TRANSACTION K WHERE K.TRANSACTION_ID = '{TRANSACTION_ID}'"
$beneIdHeader       = "bene_id_query"
$beneIdTranHeader   = "bene_id_tran_query"
$beneBankNameHeader = "bene_bank_name_query"
$beneBankNameSql    = "SELECT DB.BANK_NAME FROM DC_FUND_TRANSFER_BANK DB WHERE DB.FUND_TRANSFER_BANK_ID = (SELECT DT.BANK_ID FROM DC_TRANSACTION DT WHERE DT.TRANSACTION_ID = '{TRANSACTION_ID}')"
$beneBankTranHeader = "bene_bank_tran_query"
$beneBankTranSql    = "SELECT K.BENEFICIARY_BANK FROM DC_TRANSACTION K WHERE K.TRANSACTION_ID = '{TRANSACTION_ID}'"
$fundTransferBeneSql = "SELECT FUND_TRANSFER_BENEFICIARY_ID FROM DC_FUND_TRANSFER_BENEFICIARY K WHERE K.CUSTOMER_INFO_ID = (Select CUSTOMER_INFO_ID from DC_CUSTOMER_INFO L WHERE L.CUSTOMER_NAME = '{customer_name}' ) and K.ACCOUNT_NO = '{account_number}'"

# --- header row (row 1) --------------------------------------------------------
$ws.Range("W1").Value = $branchCodeHeader
$ws.Range("X1").Value = $toBranchHeader
$ws.Range("Y1").Value = $beneIdTranHeader
$ws.Range("Z1").Value = $beneIdHeader
$ws.Range("AA1").Value = $beneBankNameHeader
$ws.Range("AB1").Value = $beneBankTranHeader

# --- data rows (rows 2-44 all share the same query text per column) -----------
$lastRow = 44
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 23).Value = $branchCodeSql       # W
    $ws.Cells.Item($r, 24).Value = $toBranchSql         # X
    $ws.Cells.Item($r, 25).Value = $beneIdSql           # Y
    $ws.Cells.Item($r, 26).Value = $fundTransferBeneSql # Z
    $ws.Cells.Item($r, 27).Value = $beneBankNameSql     # AA
    $ws.Cells.Item($r, 28).Value = $beneBankTranSql     # AB
}

# --- column widths (best-fit widths authored alongside the new columns) -------
$ws.Columns("W").ColumnWidth = 96
$ws.Columns("X").ColumnWidth = 112.66666666666667
$ws.Columns("Y").ColumnWidth = 91.66666666666667
$ws.Columns("Z").ColumnWidth = 244.66666666666666
$ws.Columns("AA").ColumnWidth = 185.16666666666666
$ws.Columns("AB").ColumnWidth = 96.66666666666667

# --- selection / view ----------------------------------------------------------
$ws.Range("AA14").Select()
